$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# New row 11 - "GuilID" property definition (工会ID)
$ws.Range("A11").Value = "GuilID"
$ws.Range("B11").Value = "object"
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "Friend"
$ws.Range("J11").Value = "工会ID"

# Match formatting used by the rest of the table (text number format, style index 1)
$ws.Range("A11:B11").NumberFormat = "@"
$ws.Range("I11:J11").NumberFormat = "@"

# Keep selection on the sheet near the newly added row
$ws.Range("E19").Select() | Out-Null
